$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps text formatting so numeric-looking values
# like "1.001" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.666.31"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.920.84"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "239.79"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "0.4941"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "0.2997"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").Value = "0.06761"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "1.927.70"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").Value = "17.20"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "0.07343"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").Value = "5.200"
$ws.Range("E13").Value = "  +2.95%  "
$ws.Range("D14").Value = "88.67"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").Value = "0.6741"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "30.637.25"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "0.000007953"
$ws.Range("D18").Value = "13.53"
$ws.Range("E18").Value = "  +2.72%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "2.157.62"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").Value = "5.385"
$ws.Range("E21").Value = "  +11.46%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "198.37"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").Value = "6.336"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("D25").Value = "9.650"
$ws.Range("E25").Value = "  +3.02%  "
$ws.Range("D26").Value = "165.13"
$ws.Range("E26").Value = "  +6.53%  "
$ws.Range("D27").Value = "18.69"
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("D28").Value = "1.958"
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("E29").Value = "  +4.70%  "
$ws.Range("D30").Value = "4.381"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "0.09187"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").Value = "4.076"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").Value = "0.05272"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "0.7435"
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("D35").Value = "1.117"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "2.726"
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("D39").Value = "0.9282"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").Value = "2.087"
$ws.Range("E40").Value = "  -2.95%  "
$ws.Range("D41").Value = "0.4465"
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").Value = "5.972"
$ws.Range("E42").Value = "  +3.45%  "
$ws.Range("D43").Value = "71.64"
$ws.Range("E43").Value = "  +23.98%  "
$ws.Range("D44").Value = "106.28"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "1.003"
$ws.Range("D46").Value = "0.1398"
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("D47").Value = "7.655"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "9.041"
$ws.Range("E48").Value = "  +3.50%  "
$ws.Range("D49").Value = "35.09"
$ws.Range("E49").Value = "  +4.15%  "
$ws.Range("D50").Value = "0.05887"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").Value = "0.4032"
$ws.Range("E51").Value = "  +2.41%  "
